# Auto-generated: apply scheduled-runner market-data refresh to Excalibur_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 222.5
$ws.Range("I9").Value = 222.5
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 222.5
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -53.5
$ws.Range("N9").ClearContents()

# Row 43
$ws.Range("H43").Value = 1672.3636
$ws.Range("I43").Value = 1349.75
$ws.Range("J43").Value = 1856.7142
$ws.Range("K43").Value = 1349.75
$ws.Range("L43").Value = 1856.7142
$ws.Range("M43").Value = -1280.75
$ws.Range("N43").Value = -1994.7142

# Row 137
$ws.Range("H137").Value = 1277076.2
$ws.Range("I137").Value = 1069287.8
$ws.Range("J137").Value = 1385126.2
$ws.Range("K137").Value = 3207863.4
$ws.Range("L137").Value = 4155378.6
$ws.Range("M137").Value = -3205313.4
$ws.Range("N137").Value = -4160478.6

# Row 138
$ws.Range("H138").Value = 2764.0508
$ws.Range("J138").Value = 3178.027
$ws.Range("L138").Value = 9534.081
$ws.Range("N138").Value = -19814.081

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1963822.2
$ws.Range("I61").Value = 2225345.2
$ws.Range("J61").Value = 2400
$ws.Range("K61").Value = 2225345.2
$ws.Range("L61").Value = 2400
$ws.Range("M61").Value = -2225133.2
$ws.Range("N61").Value = -2824

# Row 102
$ws.Range("H102").Value = 3902.2
$ws.Range("I102").Value = 3295.25
$ws.Range("K102").Value = 3295.25
$ws.Range("M102").Value = -1673.25

# Row 110
$ws.Range("H110").Value = 2003.3334
$ws.Range("I110").Value = 2003.3334
$ws.Range("K110").Value = 2003.3334
$ws.Range("M110").Value = 41.66660000000002

# Row 132
$ws.Range("H132").Value = 2121041
$ws.Range("I132").Value = 2910806.5
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 8732419.5
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -8729889.5
$ws.Range("N132").Value = -50060

# Row 136
$ws.Range("H136").Value = 1963822.2
$ws.Range("I136").Value = 2225345.2
$ws.Range("J136").Value = 2400
$ws.Range("K136").Value = 6676035.600000001
$ws.Range("L136").Value = 7200
$ws.Range("M136").Value = -6673485.600000001
$ws.Range("N136").Value = -12300

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 42651.46
$ws.Range("I99").Value = 65061.125
$ws.Range("K99").Value = 65061.125
$ws.Range("M99").Value = -63563.125

# Row 107
$ws.Range("H107").Value = 5607.25
$ws.Range("I107").Value = 5347
$ws.Range("K107").Value = 5347
$ws.Range("M107").Value = -3427

# Row 108
$ws.Range("H108").Value = 100684
$ws.Range("J108").Value = 100684
$ws.Range("L108").Value = 100684
$ws.Range("N108").Value = -108364

$ws = $wb.Worksheets.Item("CRP")
# Row 10
$ws.Range("H10").Value = 5000
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

# Row 69
$ws.Range("H69").Value = 59498.082
$ws.Range("I69").Value = 44164.832
$ws.Range("K69").Value = 44164.832
$ws.Range("M69").Value = -43415.832

# Row 72
$ws.Range("H72").Value = 59498.082
$ws.Range("I72").Value = 44164.832
$ws.Range("K72").Value = 132494.496
$ws.Range("M72").Value = -128750.496

# Row 99
$ws.Range("H99").Value = 5071.0415
$ws.Range("I99").Value = 3909.923
$ws.Range("J99").Value = 6443.273
$ws.Range("K99").Value = 3909.923
$ws.Range("L99").Value = 6443.273
$ws.Range("M99").Value = -2411.923
$ws.Range("N99").Value = -9439.273000000001

# Row 126
$ws.Range("H126").Value = 5071.0415
$ws.Range("I126").Value = 3909.923
$ws.Range("J126").Value = 6443.273
$ws.Range("K126").Value = 11729.769
$ws.Range("L126").Value = 19329.819
$ws.Range("M126").Value = -9259.769
$ws.Range("N126").Value = -24269.819

# Row 134
$ws.Range("H134").Value = 1323.9032
$ws.Range("I134").Value = 1213.68
$ws.Range("J134").Value = 1783.1666
$ws.Range("K134").Value = 3641.04
$ws.Range("L134").Value = 5349.4998
$ws.Range("M134").Value = -1106.04
$ws.Range("N134").Value = -10419.4998

$ws = $wb.Worksheets.Item("CUL")
# Row 93
$ws.Range("H93").Value = 4122.8887
$ws.Range("J93").Value = 5351.1665
$ws.Range("L93").Value = 16053.4995
$ws.Range("N93").Value = -19797.4995

# Row 119
$ws.Range("H119").Value = 1999
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

# Row 132
$ws.Range("H132").Value = 3907.6365
$ws.Range("I132").Value = 2622.5
$ws.Range("J132").Value = 4642
$ws.Range("K132").Value = 23602.5
$ws.Range("L132").Value = 41778
$ws.Range("M132").Value = -21072.5
$ws.Range("N132").Value = -46838

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 248.16667
$ws.Range("I2").Value = 102.42857
$ws.Range("K2").Value = 102.42857
$ws.Range("M2").Value = 10.57143000000001

# Row 18
$ws.Range("H18").Value = 20006
$ws.Range("J18").Value = 20006
$ws.Range("L18").Value = 20006
$ws.Range("N18").Value = -20592

# Row 43
$ws.Range("H43").Value = 16580.334
$ws.Range("I43").Value = 3474.8
$ws.Range("J43").Value = 32962.25
$ws.Range("K43").Value = 3474.8
$ws.Range("L43").Value = 32962.25
$ws.Range("M43").Value = -3323.8
$ws.Range("N43").Value = -33264.25

# Row 46
$ws.Range("H46").Value = 29339.182
$ws.Range("I46").Value = 18624.5
$ws.Range("J46").Value = 35461.855
$ws.Range("K46").Value = 18624.5
$ws.Range("L46").Value = 35461.855
$ws.Range("M46").Value = -18468.5
$ws.Range("N46").Value = -35773.855

# Row 57
$ws.Range("H57").Value = 31999.334
$ws.Range("J57").Value = 31999.334
$ws.Range("L57").Value = 31999.334
$ws.Range("N57").Value = -33639.334

# Row 80
$ws.Range("H80").Value = 8063.273
$ws.Range("J80").Value = 9179.799999999999
$ws.Range("L80").Value = 9179.799999999999
$ws.Range("N80").Value = -11175.8

# Row 83
$ws.Range("H83").Value = 8063.273
$ws.Range("J83").Value = 9179.799999999999
$ws.Range("L83").Value = 45899
$ws.Range("N83").Value = -55883

# Row 122
$ws.Range("H122").Value = 3330.4827
$ws.Range("I122").Value = 3201.6
$ws.Range("J122").Value = 3616.889
$ws.Range("K122").Value = 9604.799999999999
$ws.Range("L122").Value = 10850.667
$ws.Range("M122").Value = -7154.799999999999
$ws.Range("N122").Value = -15750.667

# Row 132
$ws.Range("H132").Value = 862710.9399999999
$ws.Range("I132").Value = 1096914.1
$ws.Range("J132").Value = 3965.6667
$ws.Range("K132").Value = 3290742.3
$ws.Range("L132").Value = 11897.0001
$ws.Range("M132").Value = -3288212.3
$ws.Range("N132").Value = -16957.0001

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 5615.625
$ws.Range("I40").Value = 5615.625
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 5615.625
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -5479.625
$ws.Range("N40").ClearContents()

# Row 42
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

# Row 49
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0

# Row 55
$ws.Range("H55").Value = 218.88889
$ws.Range("J55").Value = 247.2
$ws.Range("L55").Value = 247.2
$ws.Range("N55").Value = -593.2

# Row 62
$ws.Range("H62").Value = 46263
$ws.Range("J62").Value = 57500
$ws.Range("L62").Value = 57500
$ws.Range("N62").Value = -58748

# Row 65
$ws.Range("H65").Value = 46263
$ws.Range("J65").Value = 57500
$ws.Range("L65").Value = 172500
$ws.Range("N65").Value = -178740

# Row 122
$ws.Range("H122").Value = 3656.4
$ws.Range("I122").Value = 3516.2727
$ws.Range("J122").Value = 4041.75
$ws.Range("K122").Value = 10548.8181
$ws.Range("L122").Value = 12125.25
$ws.Range("M122").Value = -8098.8181
$ws.Range("N122").Value = -17025.25

# Row 132
$ws.Range("H132").Value = 1283159.5
$ws.Range("I132").Value = 1442867
$ws.Range("K132").Value = 4328601
$ws.Range("M132").Value = -4326071

$ws = $wb.Worksheets.Item("WVR")
# Row 23
$ws.Range("H23").Value = 2230
$ws.Range("I23").Value = 1011.6667
$ws.Range("J23").Value = 4666.6665
$ws.Range("K23").Value = 1011.6667
$ws.Range("L23").Value = 4666.6665
$ws.Range("M23").Value = -782.6667
$ws.Range("N23").Value = -5124.6665

# Row 107
$ws.Range("H107").Value = 2070.3547
$ws.Range("I107").Value = 811.7059
$ws.Range("J107").Value = 3598.7144
$ws.Range("K107").Value = 2435.1177
$ws.Range("L107").Value = 10796.1432
$ws.Range("M107").Value = -515.1177000000002
$ws.Range("N107").Value = -14636.1432

# Row 113
$ws.Range("H113").Value = 3687.1
$ws.Range("J113").Value = 4179.7144
$ws.Range("L113").Value = 12539.1432
$ws.Range("N113").Value = -16879.1432

# Row 122
$ws.Range("H122").Value = 3766.1538
$ws.Range("I122").Value = 4077
$ws.Range("J122").Value = 3499.7144
$ws.Range("K122").Value = 12231
$ws.Range("L122").Value = 10499.1432
$ws.Range("M122").Value = -9781
$ws.Range("N122").Value = -15399.1432

# Row 132
$ws.Range("H132").Value = 8753903
$ws.Range("I132").Value = 9586418
$ws.Range("K132").Value = 28759254
$ws.Range("M132").Value = -28756724

Write-Host "Applied scheduled market-data refresh."